# Rebuild sheet1 content to match target layout (docs/assets/disciplinas/LOM3221.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing rows/content so the sheet (and shared-string table) is rebuilt cleanly
$ws.Range("A1:C25").EntireRow.Delete()

function Set-TextCell($r, $c, $text, $styleSrcR, $styleSrcC) {
    $cell = $ws.Cells.Item($r, $c)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy($cell)
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues - keep as plain text, avoid auto number/date conversion
    if ($styleSrcR) {
        $src = $ws.Cells.Item($styleSrcR, $styleSrcC)
        $src.Copy()
        $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats, restore the columns normal style
    }
    $ws.Application.CutCopyMode = $false
}

$ws.Cells.Item(1,2).Value = "Ementa atual:"
$ws.Cells.Item(1,3).Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Cells.Item(2,2).Value = "LOM3221"
$ws.Cells.Item(2,3).Value = "LOM3221"

$ws.Cells.Item(3,1).Value = "Nome:"
$ws.Cells.Item(3,2).Value = " Laboratório de Eletrônica"
$ws.Cells.Item(3,3).Value = " Laboratório de Eletrônica"

$ws.Cells.Item(4,1).Value = "Name:"
$ws.Cells.Item(4,2).Value = "Laboratory of Electronics"
$ws.Cells.Item(4,3).Value = "Laboratory of Electronics"

$ws.Cells.Item(5,1).Value = "Créditos-aula:"
Set-TextCell 5 2 "4" 3 2
Set-TextCell 5 3 "4" 3 3

$ws.Cells.Item(6,1).Value = "Créditos-trabalho"
Set-TextCell 6 2 "0" 3 2
Set-TextCell 6 3 "0" 3 3

$ws.Cells.Item(7,1).Value = "Carga horária:"
$ws.Cells.Item(7,2).Value = "60 h"
$ws.Cells.Item(7,3).Value = "60 h"

$ws.Cells.Item(8,1).Value = "Ativação:"
Set-TextCell 8 2 "01/01/2012" 3 2
Set-TextCell 8 3 "01/01/2012" 3 3

$ws.Cells.Item(9,1).Value = "Semestre ideal:"
$ws.Cells.Item(9,2).Value = "EF-5"
$ws.Cells.Item(9,3).Value = "EF-5"

$ws.Cells.Item(10,1).Value = "Objetivos:"
$ws.Cells.Item(10,2).Value = "3268262 - Carlos Renato Menegatti"
$ws.Cells.Item(10,3).Value = "3268262 - Carlos Renato Menegatti"
$ws.Rows.Item(10).RowHeight = 60

$ws.Cells.Item(11,1).Value = "Objectives:"
$ws.Rows.Item(11).RowHeight = 60

$ws.Cells.Item(12,1).Value = "Docentes responsáveis:"

$ws.Cells.Item(13,1).Value = "Programa resumido:"
$ws.Cells.Item(13,2).Value = "Semestral"
$ws.Cells.Item(13,3).Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

$ws.Cells.Item(14,1).Value = "Short syllabus:"
$ws.Rows.Item(14).RowHeight = 60

$ws.Cells.Item(15,1).Value = "Programa:"
Set-TextCell 15 2 "01/01/2012" 3 2
Set-TextCell 15 3 "01/01/2012" 3 3
$ws.Rows.Item(15).RowHeight = 120

$ws.Cells.Item(16,1).Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

$ws.Cells.Item(17,1).Value = "Avaliação:"

$ws.Cells.Item(18,1).Value = "Método:"
$ws.Cells.Item(18,2).Value = "3268262 - Carlos Renato Menegatti"
$ws.Cells.Item(18,3).Value = "3268262 - Carlos Renato Menegatti"
$ws.Rows.Item(18).RowHeight = 60

$ws.Cells.Item(19,1).Value = "Critério:"
$ws.Cells.Item(19,2).Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento, desenvolvimento e apresentação de projetos de circuitos e realização de testes sobre o experimento em estudo."
$ws.Cells.Item(19,3).Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento, desenvolvimento e apresentação de projetos de circuitos e realização de testes sobre o experimento em estudo."
$ws.Rows.Item(19).RowHeight = 60

$ws.Cells.Item(20,1).Value = "Norma de recuperação:"
$ws.Cells.Item(20,2).Value = "Média aritmética de prova escrita, testes, trabalhos e relatórios: PE, TS e TR. Conceito Final = (PE  + TS + TR)/3"
$ws.Cells.Item(20,3).Value = "Média aritmética de prova escrita, testes, trabalhos e relatórios: PE, TS e TR. Conceito Final = (PE  + TS + TR)/3"
$ws.Rows.Item(20).RowHeight = 60

$ws.Cells.Item(21,1).Value = "Bibliografia:"
$ws.Cells.Item(21,2).Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Cells.Item(21,3).Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Rows.Item(21).RowHeight = 120

$ws.Cells.Item(22,1).Value = "Requisitos:"

$ws.Cells.Item(23,2).Value = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$ws.Cells.Item(23,3).Value = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

$ws.Cells.Item(24,2).Value = "LOM3206 -  Eletrônica  (Indicação de Conjunto)`n"
$ws.Cells.Item(24,3).Value = "LOM3206 -  Eletrônica  (Indicação de Conjunto)`n"
$ws.Rows.Item(24).RowHeight = 30

